$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.502.30"
$ws.Range("E2").Value = "  +5.13%  "

# Row 3
$ws.Range("D3").Value = "1.602.06"

# Row 4
$ws.Range("E4").Value = "  -0.17%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.10"
$ws.Range("E5").Value = "  +2.46%  "

# Row 6
$ws.Range("E6").Value = "  +1.87%  "

# Row 7
$ws.Range("E7").Value = "  -0.15%  "

# Row 8
$ws.Range("E8").Value = "  +9.36%  "

# Row 9
$ws.Range("E9").Value = "  +1.76%  "

# Row 10
$ws.Range("E10").Value = "  +1.19%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0891"
$ws.Range("E11").Value = "  +2.49%  "

# Row 12
$ws.Range("D12").Value = "1.831.46"
$ws.Range("E12").Value = "  +2.83%  "

# Row 13
$ws.Range("D13").Value = "1.604.46"
$ws.Range("E13").Value = "  +3.25%  "

# Row 14
$ws.Range("E14").Value = "  +1.00%  "

# Row 15
$ws.Range("E15").Value = "  +3.65%  "

# Row 16
$ws.Range("D16").Value = "28.512.25"
$ws.Range("E16").Value = "  +5.21%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.42"
$ws.Range("E17").Value = "  +2.81%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "232.58"
$ws.Range("E18").Value = "  +7.79%  "

# Row 19
$ws.Range("E19").Value = "  +1.66%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0712"
$ws.Range("E20").Value = "  +1.70%  "

# Row 21
$ws.Range("E21").Value = "  -0.23%  "

# Row 22
$ws.Range("E22").Value = "  +0.76%  "

# Row 23
$ws.Range("E23").Value = "  +2.73%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.96"
$ws.Range("E24").Value = "  +1.63%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.41"
$ws.Range("E25").Value = "  -0.32%  "

# Row 26
$ws.Range("E26").Value = "  +2.36%  "

# Row 27
$ws.Range("E27").Value = "  +0.46%  "

# Row 28
$ws.Range("E28").Value = "  +1.30%  "

# Row 29
$ws.Range("E29").Value = "  -0.10%  "

# Row 30
$ws.Range("E30").Value = "  +1.33%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0476"
$ws.Range("E31").Value = "  +1.57%  "

# Row 32
$ws.Range("E32").Value = "  +1.02%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.17"
$ws.Range("E33").Value = "  +0.56%  "

# Row 34
$ws.Range("D34").Value = "1.425.46"
$ws.Range("E34").Value = "  -0.60%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.61"
$ws.Range("E35").Value = "  -0.09%  "

# Row 36
$ws.Range("E36").Value = "  -4.08%  "

# Row 37
$ws.Range("E37").Value = "  -0.09%  "

# Row 38
$ws.Range("E38").Value = "  +1.33%  "

# Row 39 (content swapped with what used to be row 40, with refreshed D/E)
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.53"
$ws.Range("E39").Value = "  +8.47%  "

# Row 40 (content swapped with what used to be row 39, with refreshed D/E)
$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.546"
$ws.Range("E40").Value = "  +3.07%  "

# Row 41
$ws.Range("E41").Value = "  +2.42%  "

# Row 42
$ws.Range("E42").Value = "  -2.64%  "

# Row 44 (content swapped with what used to be row 45, with refreshed D/E)
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.985"
$ws.Range("E44").Value = "  -1.18%  "

# Row 45 (content swapped with what used to be row 44, with refreshed D/E)
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.84"
$ws.Range("E45").Value = "  +6.77%  "

# Row 46
$ws.Range("E46").Value = "  +1.41%  "

# Row 47
$ws.Range("D47").Value = "1.742.01"
$ws.Range("E47").Value = "  +2.90%  "

# Row 48 (content swapped with what used to be row 49, with refreshed D/E)
$ws.Range("B48").Value = "mCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.14"
$ws.Range("E48").Value = "  +0.36%  "

# Row 49 (content swapped with what used to be row 48, with refreshed D/E)
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "87.57"
$ws.Range("E49").Value = "  +2.57%  "

# Row 50
$ws.Range("D50").Value = "0.0₆0108"
$ws.Range("E50").Value = "  +9.23%  "

# Row 51
$ws.Range("E51").Value = "  +0.63%  "
